$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.002.94"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "3.335.76"
$ws.Range("E3").Value = "  +1.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.50"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.35"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("D9").Value = "3.331.89"
$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("E10").Value = "  +4.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.84"
$ws.Range("E12").Value = "  +3.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "684.82"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").Value = "3.881.85"
$ws.Range("E15").Value = "  +2.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("D17").Value = "68.073.12"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").Value = "3.347.34"
$ws.Range("E18").Value = "  +2.24%  "

$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.43"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("E21").Value = "  +3.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("E22").Value = "  +1.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.12"
$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  +4.30%  "

$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("E28").Value = "  +2.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.08"
$ws.Range("E31").Value = "  +4.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "574.57"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.03"
$ws.Range("E33").Value = "  +2.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  +1.99%  "

$ws.Range("D35").Value = "3.713.87"
$ws.Range("E35").Value = "  -4.02%  "

$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.08"
$ws.Range("E37").Value = "  +3.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.62"
$ws.Range("E39").Value = "  +9.35%  "

$ws.Range("E40").Value = "  +2.45%  "

$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("E42").Value = "  +5.76%  "

$ws.Range("D43").Value = "0.0₃0677"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.336"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("E47").Value = "  +6.17%  "

$ws.Range("E48").Value = "  +1.37%  "

$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.62"
$ws.Range("E51").Value = "  +0.31%  "
